$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header F1 from "eMail" to "GebDat"
$ws.Range("F1").Value = "GebDat"

# Fill birth dates for rows 2-7 (as text, matching source diff which stores them as shared strings, not numeric dates)
$ws.Range("F2").Value = "30.05.1982"
$ws.Range("F3").Value = "11.09.1998"
$ws.Range("F4").Value = "30.01.1986"
$ws.Range("F5").Value = "02.01.1990"
$ws.Range("F6").Value = "22.10.1990"
$ws.Range("F7").Value = "17.03.1988"
